$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9450
$ws.Range("E2").Value = 428
$ws.Range("F2").Value = 428
$ws.Range("G2").Value = 374
$ws.Range("H2").Value = 279
$ws.Range("I2").Value = 279
$ws.Range("K2").Value = 5638
$ws.Range("L2").Value = 3849
$ws.Range("M2").Value = 1789
$ws.Range("N2").Value = 1789
$ws.Range("P2").Value = 172
$ws.Range("Q2").Value = 211
$ws.Range("R2").Value = -169
$ws.Range("S2").Value = -49
$ws.Range("T2").Value = 137
$ws.Range("U2").Value = 75
$ws.Range("V2").Value = 2475
$ws.Range("W2").Value = 4.53
$ws.Range("X2").Value = 2.96
$ws.Range("Y2").Value = 16.99
$ws.Range("Z2").Value = 5.03
$ws.Range("AA2").Value = 215.1
$ws.Range("AB2").Value = 918.66
$ws.Range("AC2").Value = 875
$ws.Range("AD2").Value = 8.16
$ws.Range("AE2").Value = 5609
$ws.Range("AF2").Value = 1.27
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.4
$ws.Range("AI2").Value = 11.42
$ws.Range("AJ2").Value = 31900000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 8677
$ws.Range("E3").Value = 607
$ws.Range("F3").Value = 607
$ws.Range("G3").Value = 602
$ws.Range("H3").Value = 449
$ws.Range("I3").Value = 449
$ws.Range("K3").Value = 5606
$ws.Range("L3").Value = 3420
$ws.Range("M3").Value = 2187
$ws.Range("N3").Value = 2187
$ws.Range("P3").Value = 172
$ws.Range("Q3").Value = 622
$ws.Range("R3").Value = -45
$ws.Range("S3").Value = -382
$ws.Range("T3").Value = 47
$ws.Range("U3").Value = 575
$ws.Range("V3").Value = 2152
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 5.17
$ws.Range("Y3").Value = 22.56
$ws.Range("Z3").Value = 7.98
$ws.Range("AA3").Value = 156.38
$ws.Range("AB3").Value = 1153.33
$ws.Range("AC3").Value = 1406
$ws.Range("AD3").Value = 7.82
$ws.Range("AE3").Value = 6855
$ws.Range("AF3").Value = 1.6
$ws.Range("AG3").Value = 130
$ws.Range("AH3").Value = 1.18
$ws.Range("AI3").Value = 9.24
$ws.Range("AJ3").Value = 31900000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 8612
$ws.Range("E4").Value = 689
$ws.Range("F4").Value = 689
$ws.Range("G4").Value = 655
$ws.Range("H4").Value = 511
$ws.Range("I4").Value = 511
$ws.Range("K4").Value = 5772
$ws.Range("L4").Value = 3121
$ws.Range("M4").Value = 2650
$ws.Range("N4").Value = 2650
$ws.Range("P4").Value = 172
$ws.Range("Q4").Value = 496
$ws.Range("R4").Value = -140
$ws.Range("S4").Value = -384
$ws.Range("T4").Value = 94
$ws.Range("U4").Value = 402
$ws.Range("V4").Value = 1835
$ws.Range("W4").Value = 8
$ws.Range("X4").Value = 5.94
$ws.Range("Y4").Value = 21.14
$ws.Range("Z4").Value = 8.99
$ws.Range("AA4").Value = 117.77
$ws.Range("AB4").Value = 1431.39
$ws.Range("AC4").Value = 1603
$ws.Range("AD4").Value = 5.84
$ws.Range("AE4").Value = 8388
$ws.Range("AF4").Value = 1.12
$ws.Range("AG4").Value = 140
$ws.Range("AH4").Value = 1.5
$ws.Range("AI4").Value = 8.67
$ws.Range("AJ4").Value = 31900000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 9144
$ws.Range("E5").Value = 415
$ws.Range("F5").Value = 415
$ws.Range("G5").Value = 433
$ws.Range("H5").Value = 310
$ws.Range("I5").Value = 310
$ws.Range("K5").Value = 6127
$ws.Range("L5").Value = 3329
$ws.Range("M5").Value = 2798
$ws.Range("N5").Value = 2798
$ws.Range("P5").Value = 172
$ws.Range("Q5").Value = 77
$ws.Range("R5").Value = -208
$ws.Range("S5").Value = 54
$ws.Range("T5").Value = 210
$ws.Range("U5").Value = -133
$ws.Range("V5").Value = 1944
$ws.Range("W5").Value = 4.54
$ws.Range("X5").Value = 3.4
$ws.Range("Y5").Value = 11.39
$ws.Range("Z5").Value = 5.22
$ws.Range("AA5").Value = 118.97
$ws.Range("AB5").Value = 1589.46
$ws.Range("AC5").Value = 973
$ws.Range("AD5").Value = 6.99
$ws.Range("AE5").Value = 9056
$ws.Range("AF5").Value = 0.75
$ws.Range("AG5").Value = 120
$ws.Range("AH5").Value = 1.76
$ws.Range("AI5").Value = 11.95
$ws.Range("AJ5").Value = 31900000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 9108
$ws.Range("E6").Value = 144
$ws.Range("F6").Value = 144
$ws.Range("G6").Value = 123
$ws.Range("H6").Value = 85
$ws.Range("I6").Value = 85
$ws.Range("K6").Value = 6004
$ws.Range("L6").Value = 3181
$ws.Range("M6").Value = 2823
$ws.Range("N6").Value = 2823
$ws.Range("P6").Value = 172
$ws.Range("Q6").Value = 82
$ws.Range("R6").Value = -78
$ws.Range("S6").Value = -59
$ws.Range("T6").Value = 79
$ws.Range("U6").Value = 3
$ws.Range("V6").Value = 1938
$ws.Range("W6").Value = 1.58
$ws.Range("X6").Value = 0.93
$ws.Range("Y6").Value = 3.02
$ws.Range("Z6").Value = 1.4
$ws.Range("AA6").Value = 112.69
$ws.Range("AB6").Value = 1617.1
$ws.Range("AC6").Value = 266
$ws.Range("AD6").Value = 16.13
$ws.Range("AE6").Value = 9135
$ws.Range("AF6").Value = 0.47
$ws.Range("AG6").Value = 60
$ws.Range("AH6").Value = 1.4
$ws.Range("AI6").Value = 21.82
$ws.Range("AJ6").Value = 31900000

# Rows 7-9: clear all data columns (D:AJ), keep A/B/C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
